$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Update the "repaymentstrategy" value scenario from "RBI (India)" to the
# new periodic/upfront related scenario text.
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Move the active selection onto the cell that was just edited.
$ws.Activate()
$ws.Range("B17").Select()
